# Add pricing_origination_charges, pricing_broker_fees, and pricing_initial_charges
# fields to the "invalid" sheet (columns U, V, W), matching the layout/format of the
# existing pricing_* columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid")
$ws.Activate()

# ---- Header row (row 1) ----
$ws.Range("U1").Value = "pricing_origination_charges"
$ws.Range("V1").Value = "pricing_broker_fees"
$ws.Range("W1").Value = "pricing_initial_charges"
$ws.Range("U1:W1").WrapText = $true

# Header row needs to grow to fit the new wrapped header text.
$ws.Rows.Item(1).RowHeight = 51

# ---- Data rows (rows 2-11) ----
# pricing_origination_charges (U), pricing_broker_fees (V), pricing_initial_charges (W)
$data = @(
  @(2,  1,    6, 1),
  @(3,  2,    5, 2),
  @(4,  "a",  "c", "s"),
  @(5,  4,    5, 2),
  @(6,  5,    5, 2),
  @(7,  6,    5, 2),
  @(8,  "   ",5, 22),
  @(9,  8,    5, 2),
  @(10, 8,    5, 2),
  @(11, 10,   5, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("U$r").Value = $row[1]
    $ws.Range("V$r").Value = $row[2]
    $ws.Range("W$r").Value = $row[3]
}

$ws.Range("U2:W11").WrapText = $true

# ---- Column widths ----
# P and Q widen to fit the new long header labels; R/S go back to the sheet default
# width; T (pricing_prepenalty_exists, pushed out by the new columns) gets its own
# custom width; everything from U onward stays at the default width.
$ws.Columns.Item(16).ColumnWidth = 18.75  # P
$ws.Columns.Item(17).ColumnWidth = 18.75  # Q
$ws.Columns.Item(18).ColumnWidth = 9.92   # R -> back to default-ish width
$ws.Columns.Item(20).ColumnWidth = 20.67  # T

# ---- View/selection bookkeeping ----
$ws.Range("T1").Select()
